$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 2-21 per repull/push of data and mean calculation.
$newF = @{
    2  = 3
    3  = -2
    4  = -1
    5  = 1
    6  = 1
    7  = 4
    8  = -3
    9  = 0
    10 = -1
    11 = 6
    12 = 0
    13 = -2
    14 = -3
    15 = 3
    16 = -3
    17 = -1
    18 = 3
    19 = 1
    20 = 0
    21 = -4
}

foreach ($row in $newF.Keys) {
    $ws.Range("F$row").Value = $newF[$row]
}
